# Generate Report for Handoff
#
# c152a9e8-e850-4b25-b9c4-efdc5c5dd72e is now "Ready for handoff" and moves
# ahead of deacab6e-5088-41e0-b5b3-4ecc97525912 in the status report (rows 7
# and 8 swap their content across the Overview/zh-cn/de-de sheets), and
# deacab6e's latest handoff timestamps are refreshed to reflect a new
# handoff event.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay {
    param($ws, [string]$addr, [string]$text)
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.TextToDisplay = $text
        }
    }
}

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A7").Value = "c152a9e8-e850-4b25-b9c4-efdc5c5dd72e.md"
$ov.Range("B7").Value = "Ready for handoff"
$ov.Range("C7").Value = "Ready for handoff"
$ov.Range("D7").Value = "2016-03-31 06:47:52"

$ov.Range("A8").Value = "deacab6e-5088-41e0-b5b3-4ecc97525912.md"
$ov.Range("B8").Value = "Ready for handoff"
$ov.Range("C8").Value = "Ready for handoff"
$ov.Range("D8").Value = "2016-03-31 06:54:03"

Set-HyperlinkDisplay $ov '$A$7' "c152a9e8-e850-4b25-b9c4-efdc5c5dd72e.md"
Set-HyperlinkDisplay $ov '$A$8' "deacab6e-5088-41e0-b5b3-4ecc97525912.md"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A7").Value = "c152a9e8-e850-4b25-b9c4-efdc5c5dd72e.md"
$zh.Range("C7").Value = "Ready for handoff"
$zh.Range("D7").Value = "c152a9e8-e850-4b25-b9c4-efdc5c5dd72e.4660fb0f7c81c1dad25cb162ae5b5d57537b745a.zh-cn.xlf"
$zh.Range("E7").Value = "2016-03-31 06:47:42"

$zh.Range("A8").Value = "deacab6e-5088-41e0-b5b3-4ecc97525912.md"
$zh.Range("C8").Value = "Ready for handoff"
$zh.Range("D8").Value = "deacab6e-5088-41e0-b5b3-4ecc97525912.d539915c154dd968b6735e46ff6a28bcd3d844ff.zh-cn.xlf"
$zh.Range("E8").Value = "2016-03-31 06:53:51"

Set-HyperlinkDisplay $zh '$A$7' "c152a9e8-e850-4b25-b9c4-efdc5c5dd72e.md"
Set-HyperlinkDisplay $zh '$D$7' "c152a9e8-e850-4b25-b9c4-efdc5c5dd72e.4660fb0f7c81c1dad25cb162ae5b5d57537b745a.zh-cn.xlf"
Set-HyperlinkDisplay $zh '$A$8' "deacab6e-5088-41e0-b5b3-4ecc97525912.md"
Set-HyperlinkDisplay $zh '$D$8' "deacab6e-5088-41e0-b5b3-4ecc97525912.d539915c154dd968b6735e46ff6a28bcd3d844ff.zh-cn.xlf"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A7").Value = "c152a9e8-e850-4b25-b9c4-efdc5c5dd72e.md"
$de.Range("C7").Value = "Ready for handoff"
$de.Range("D7").Value = "c152a9e8-e850-4b25-b9c4-efdc5c5dd72e.4660fb0f7c81c1dad25cb162ae5b5d57537b745a.de-de.xlf"
$de.Range("E7").Value = "2016-03-31 06:47:52"

$de.Range("A8").Value = "deacab6e-5088-41e0-b5b3-4ecc97525912.md"
$de.Range("C8").Value = "Ready for handoff"
$de.Range("D8").Value = "deacab6e-5088-41e0-b5b3-4ecc97525912.d539915c154dd968b6735e46ff6a28bcd3d844ff.de-de.xlf"
$de.Range("E8").Value = "2016-03-31 06:54:03"

Set-HyperlinkDisplay $de '$A$7' "c152a9e8-e850-4b25-b9c4-efdc5c5dd72e.md"
Set-HyperlinkDisplay $de '$D$7' "c152a9e8-e850-4b25-b9c4-efdc5c5dd72e.4660fb0f7c81c1dad25cb162ae5b5d57537b745a.de-de.xlf"
Set-HyperlinkDisplay $de '$A$8' "deacab6e-5088-41e0-b5b3-4ecc97525912.md"
Set-HyperlinkDisplay $de '$D$8' "deacab6e-5088-41e0-b5b3-4ecc97525912.d539915c154dd968b6735e46ff6a28bcd3d844ff.de-de.xlf"
